# Add data for 2021-09-30
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2021-09-22"

# Update the September label in column A (row 10)
$ws.Range("A10").Value = "September (through 09-22)"

# Update September row (row 10) figures for 2017-2021 (columns D-H)
$ws.Range("D10").Value = 51
$ws.Range("E10").Value = 42
$ws.Range("F10").Value = 55
$ws.Range("G10").Value = 84
$ws.Range("H10").Value = 133

# Update Total row (row 11) figures for 2017-2021 (columns D-H)
$ws.Range("D11").Value = 602
$ws.Range("E11").Value = 532
$ws.Range("F11").Value = 404
$ws.Range("G11").Value = 868
$ws.Range("H11").Value = 1203
